$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.3630829741045147
$ws.Range("C2").Value = 0.07436217548442414
$ws.Range("D2").Value = 0.07896829115506421
$ws.Range("E2").Value = 0.417430411203398
$ws.Range("G2").Value = 0.4730071219003236
$ws.Range("H2").Value = 0.615445443545795
$ws.Range("K2").Value = 0.3611517495741055
$ws.Range("O2").Value = 2.135659816244896
$ws.Range("B3").Value = 0.3199292133546692
$ws.Range("C3").Value = 0.06872356251116685
$ws.Range("D3").Value = 0.07157345869859455
$ws.Range("E3").Value = 0.3642290642159907
$ws.Range("G3").Value = 0.4725992718885834
$ws.Range("H3").Value = 0.6195510693863042
$ws.Range("K3").Value = 0.3150789838560968
$ws.Range("O3").Value = 2.143091536191079
$ws.Range("B4").Value = 0.2934177868623919
$ws.Range("C4").Value = 0.06523826926432719
$ws.Range("D4").Value = 0.06706671304205258
$ws.Range("E4").Value = 0.3316472113900346
$ws.Range("G4").Value = 0.4727501874685984
$ws.Range("H4").Value = 0.6223997372522518
$ws.Range("K4").Value = 0.2867111767695292
$ws.Range("O4").Value = 2.14918138829654
$ws.Range("B5").Value = 0.2826108945052965
$ws.Range("C5").Value = 0.06381222610910697
$ws.Range("D5").Value = 0.06523865918410365
$ws.Range("E5").Value = 0.3183893663237711
$ws.Range("G5").Value = 0.4729122734776752
$ws.Range("H5").Value = 0.6236429589054353
$ws.Range("K5").Value = 0.275131698461621
$ws.Range("O5").Value = 2.152046125000766
$ws.Range("B6").Value = 0.2808162344049947
$ws.Range("C6").Value = 0.06357508694452463
$ws.Range("D6").Value = 0.06493562464811475
$ws.Range("E6").Value = 0.3161890440252222
$ws.Range("G6").Value = 0.4729452531881009
$ws.Range("H6").Value = 0.6238543678727595
$ws.Range("K6").Value = 0.2732077795866417
$ws.Range("O6").Value = 2.152544924433158
$ws.Range("B7").Value = 0.2932720539607487
$ws.Range("C7").Value = 0.06521906037207259
$ws.Range("D7").Value = 0.0670420249207524
$ws.Range("E7").Value = 0.3314683344940335
$ws.Range("G7").Value = 0.4727519665891009
$ws.Range("H7").Value = 0.622416170359493
$ws.Range("K7").Value = 0.2865550896443096
$ws.Range("O7").Value = 2.149218473215967
$ws.Range("B8").Value = 0.3482069005472965
$ws.Range("C8").Value = 0.07242282701233194
$ws.Range("D8").Value = 0.07641155650675557
$ws.Range("E8").Value = 0.3990681778041392
$ws.Range("G8").Value = 0.4727829957343914
$ws.Range("H8").Value = 0.616792991875343
$ws.Range("K8").Value = 0.3452825197934715
$ws.Range("O8").Value = 2.137904886765
$ws.Range("B9").Value = 0.4558013648214683
$ws.Range("C9").Value = 0.08636336574257086
$ws.Range("D9").Value = 0.09505346380970536
$ws.Range("E9").Value = 0.5323758305734856
$ws.Range("G9").Value = 0.4760440772631114
$ws.Range("H9").Value = 0.6083699095820094
$ws.Range("K9").Value = 0.4598048714981928
$ws.Range("O9").Value = 2.127873202394539
$ws.Range("B10").Value = 0.5347582878139292
$ws.Range("C10").Value = 0.09649005007750588
$ws.Range("D10").Value = 0.1089157844950819
$ws.Range("E10").Value = 0.6308879885692988
$ws.Range("G10").Value = 0.4804139600839932
$ws.Range("H10").Value = 0.6037735254232359
$ws.Range("K10").Value = 0.5435405329928358
$ws.Range("O10").Value = 2.127971385348189
$ws.Range("B11").Value = 0.5706557351390416
$ws.Range("C11").Value = 0.1010715211841386
$ws.Range("D11").Value = 0.1152588159029051
$ws.Range("E11").Value = 0.6758536141595357
$ws.Range("G11").Value = 0.4828355917475875
$ws.Range("H11").Value = 0.6020292518516044
$ws.Range("K11").Value = 0.5815442734934493
$ws.Range("O11").Value = 2.129650699848781
$ws.Range("B12").Value = 0.5842458832802322
$ws.Range("C12").Value = 0.102802729498876
$ws.Range("D12").Value = 0.1176660942419971
$ws.Range("E12").Value = 0.6929046983389213
$ws.Range("G12").Value = 0.4838153600685331
$ws.Range("H12").Value = 0.6014186704977647
$ws.Range("K12").Value = 0.5959222861862656
$ws.Range("O12").Value = 2.130522673439941
$ws.Range("B13").Value = 0.58131916237366
$ws.Range("C13").Value = 0.1024300482763891
$ws.Range("D13").Value = 0.1171474077689538
$ws.Range("E13").Value = 0.6892313659905085
$ws.Range("G13").Value = 0.4836015526562534
$ws.Range("H13").Value = 0.6015479477022723
$ws.Range("K13").Value = 0.5928263182657929
$ws.Range("O13").Value = 2.130324363934818
$ws.Range("B14").Value = 0.5717738764899707
$ws.Range("C14").Value = 0.1012140231957375
$ws.Range("D14").Value = 0.1154567577755472
$ws.Range("E14").Value = 0.6772559346362073
$ws.Range("G14").Value = 0.4829149381373554
$ws.Range("H14").Value = 0.6019780176211356
$ws.Range("K14").Value = 0.5827274283395525
$ws.Range("O14").Value = 2.129717701078221
$ws.Range("B15").Value = 0.5659266511758574
$ws.Range("C15").Value = 0.1004686894011257
$ws.Range("D15").Value = 0.1144218770734682
$ws.Range("E15").Value = 0.6699237591827654
$ws.Range("G15").Value = 0.4825025498605839
$ws.Range("H15").Value = 0.6022479538493428
$ws.Range("K15").Value = 0.5765398350100668
$ws.Range("O15").Value = 2.129376873224231
$ws.Range("B16").Value = 0.5324118512309326
$ws.Range("C16").Value = 0.09619012686701467
$ws.Range("D16").Value = 0.1085019960520128
$ws.Range("E16").Value = 0.6279525841367786
$ws.Range("G16").Value = 0.4802644610712292
$ws.Range("H16").Value = 0.6038945008637739
$ws.Range("K16").Value = 0.5410550831640819
$ws.Range("O16").Value = 2.127894616357821
$ws.Range("B17").Value = 0.511846014943302
$ws.Range("C17").Value = 0.09355886045817385
$ws.Range("D17").Value = 0.1048798087435756
$ws.Range("E17").Value = 0.6022447701547549
$ws.Range("G17").Value = 0.4790028274480136
$ws.Range("H17").Value = 0.6049934519980127
$ws.Range("K17").Value = 0.5192634184252825
$ws.Range("O17").Value = 2.127404678651146
$ws.Range("B18").Value = 0.5000151940395767
$ws.Range("C18").Value = 0.09204305703411819
$ws.Range("D18").Value = 0.1027999048679504
$ws.Range("E18").Value = 0.587472513699268
$ws.Range("G18").Value = 0.4783179732380916
$ws.Range("H18").Value = 0.6056581620089361
$ws.Range("K18").Value = 0.5067211740253583
$ws.Range("O18").Value = 2.127276680936063
$ws.Range("B19").Value = 0.4960091709081382
$ws.Range("C19").Value = 0.09152942739541459
$ws.Range("D19").Value = 0.1020962835524131
$ws.Range("E19").Value = 0.5824732676390596
$ws.Range("G19").Value = 0.4780930900571434
$ws.Range("H19").Value = 0.6058888215946467
$ws.Range("K19").Value = 0.5024731836701903
$ws.Range("O19").Value = 2.127259728469824
$ws.Range("B20").Value = 0.5140354831351885
$ws.Range("C20").Value = 0.09383920891168884
$ws.Range("D20").Value = 0.1052650364131011
$ws.Range("E20").Value = 0.6049799289325648
$ws.Range("G20").Value = 0.4791329047051676
$ws.Range("H20").Value = 0.6048730898257872
$ws.Range("K20").Value = 0.5215840365069937
$ws.Range("O20").Value = 2.127440908283091
$ws.Range("B21").Value = 0.5745776555123143
$ws.Range("C21").Value = 0.1015713000860927
$ws.Range("D21").Value = 0.1159531986243394
$ws.Range("E21").Value = 0.6807727553818381
$ws.Range("G21").Value = 0.483114907475553
$ws.Range("H21").Value = 0.6018503396572612
$ws.Range("K21").Value = 0.5856940789072951
$ws.Range("O21").Value = 2.129889478310474
$ws.Range("B22").Value = 0.6141252102490853
$ws.Range("C22").Value = 0.1066030972092591
$ws.Range("D22").Value = 0.1229694767439042
$ws.Range("E22").Value = 0.7304459901275067
$ws.Range("G22").Value = 0.48608330043065
$ws.Range("H22").Value = 0.6001659122755143
$ws.Range("K22").Value = 0.6275166423982341
$ws.Range("O22").Value = 2.132866124336061
$ws.Range("B23").Value = 0.5930199517223116
$ws.Range("C23").Value = 0.1039195292049726
$ws.Range("D23").Value = 0.1192219286517258
$ws.Range("E23").Value = 0.7039212260471004
$ws.Range("G23").Value = 0.4844654067442065
$ws.Range("H23").Value = 0.6010382558783789
$ws.Range("K23").Value = 0.60520237425294
$ws.Range("O23").Value = 2.131151166190875
$ws.Range("B24").Value = 0.5130456471055993
$ws.Range("C24").Value = 0.09371247288144957
$ws.Range("D24").Value = 0.1050908670735708
$ws.Range("E24").Value = 0.6037433404252397
$ws.Range("G24").Value = 0.4790739707222684
$ws.Range("H24").Value = 0.604927403057971
$ws.Range("K24").Value = 0.5205349284145484
$ws.Range("O24").Value = 2.127424050237522
$ws.Range("B25").Value = 0.4267096537429893
$ws.Range("C25").Value = 0.08261220862299012
$ws.Range("D25").Value = 0.08998134715716333
$ws.Range("E25").Value = 0.4962208113064719
$ws.Range("G25").Value = 0.4748168836042339
$ws.Range("H25").Value = 0.610369331458287
$ws.Range("K25").Value = 0.4288935849608606
$ws.Range("O25").Value = 2.129279764739834
